# Weekly data refresh for the "Perejil / Vega Modelo de Temuco" price log.
#
# The sheet holds one row per weekly price record (rows 2..249), all sharing
# the same fixed metadata (Mercado/Region/Categoria/etc.) but each carrying
# its own Fecha (D), Volumen (J), Precio min/max/prom (K/L/M), Unidad (N),
# Origen (O), Precio $/Kg (P) and Kg o Unidades (Q).
#
# A brand-new weekly record is being inserted at row 145; every existing
# record from row 145 down to row 249 is pushed one row further down (row
# 249's record lands in a newly created row 250), and row 145 receives the
# new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 145
$lastRow  = 249
$newRow   = $lastRow + 1

# 1) Preserve the record currently in the last row (249) so it can be moved
#    into the brand-new row (250) before anything gets overwritten.
$lastDate  = $ws.Range("D$lastRow").Value2
$lastBlock = $ws.Range("J$lastRow`:Q$lastRow").Value2
$lastFmt   = $ws.Range("D$lastRow").NumberFormat

# 2) Push every record down by one row, working from the bottom up so that
#    each source row is read before it gets overwritten.
for ($n = $lastRow; $n -gt $firstRow; $n--) {
    $prev = $n - 1
    $ws.Range("D$n").Value2        = $ws.Range("D$prev").Value2
    $ws.Range("J$n`:Q$n").Value2   = $ws.Range("J$prev`:Q$prev").Value2
}

# 3) Write the new weekly record into the now-vacated first row.
$ws.Range("D$firstRow").Value2 = 44574
$ws.Range("J$firstRow").Value2 = 125

# 4) Create the new last row, carrying over the fixed metadata from the row
#    above it plus the record that was displaced out of row 249.
$ws.Range("A$newRow").Value2 = $ws.Range("A$lastRow").Value2
$ws.Range("B$newRow").Value2 = $ws.Range("B$lastRow").Value2
$ws.Range("C$newRow").Value2 = $ws.Range("C$lastRow").Value2
$ws.Range("D$newRow").Value2 = $lastDate
$ws.Range("D$newRow").NumberFormat = $lastFmt
$ws.Range("E$newRow").Value2 = $ws.Range("E$lastRow").Value2
$ws.Range("F$newRow").Value2 = $ws.Range("F$lastRow").Value2
$ws.Range("G$newRow").Value2 = $ws.Range("G$lastRow").Value2
$ws.Range("H$newRow").Value2 = $ws.Range("H$lastRow").Value2
$ws.Range("I$newRow").Value2 = $ws.Range("I$lastRow").Value2
$ws.Range("J$newRow`:Q$newRow").Value2 = $lastBlock
$ws.Range("R$newRow").Value2 = $ws.Range("R$lastRow").Value2
